$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "mangler kun konkl" -----------------------------------------------
# C14 ("L04 - Pipelines" / Qd) used to read
# "mangler, ændre nogle parametre blabla" on a red background.
# That remark is resolved now: wipe the text and mark it done (green).
$ws.Range("C14").ClearContents()
$ws.Range("C14").Interior.Color = 5287936

# C17 ("L05 - Train linear regression" / Qa) still says "Mangler alt",
# but is no longer considered missing -> flip it from red to green.
$ws.Range("C17").Interior.Color = 5287936

# D20 ("L06 - ANN" header row) used to hold "mangler intro" on an
# orange background. The intro note is gone now - clear it and give it
# a plain flat (white/theme) fill instead of leaving it colored.
$ws.Range("D20").ClearContents()
$ws.Range("D20").Interior.ThemeColor = 2
$ws.Range("D20").Interior.TintAndShade = 0

# --- "... og ret billede" ------------------------------------------------
# The remaining L06 - ANN remarks (Qa-Qf) are all resolved too; flip
# every one of them from orange/red over to green, keeping their text.
$ws.Range("C21").Interior.Color = 5287936
$ws.Range("C22").Interior.Color = 5287936
$ws.Range("C23").Interior.Color = 5287936
$ws.Range("C24").Interior.Color = 5287936
$ws.Range("C25").Interior.Color = 5287936
$ws.Range("C26").Interior.Color = 5287936

# --- update the saved view/selection state ------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
